# Regenerate orders with updated distance/size codes.
# The workbook encodes trial conditions as text tokens such as
# "Face08_D64_S30", "Fixation_D64_l.png", "D64", "S30", etc.
# This edit renames the distance codes D64/D51/D80 -> D69/D55/D86
# and the size code S30 -> S31 everywhere they occur on the sheet
# (condition names, filenames, and the Distance/Size lookup lists).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$used = $ws.UsedRange

# Distance code renames (order-independent; codes are distinct 2-digit numbers
# that don't collide with each other or with the size codes).
$used.Replace("D64", "D69") | Out-Null
$used.Replace("D51", "D55") | Out-Null
$used.Replace("D80", "D86") | Out-Null

# Size code rename.
$used.Replace("S30", "S31") | Out-Null
